$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Experimental Plan")

# Insert a new blank row above row 10 (shifts the existing "Model: Random
# Forest" ... "Finalise Experimental Report and Documentation" rows down by
# one, B10:D19 -> B11:D20).
$ws.Rows.Item(10).Insert()

# The newly inserted row has no formatting; clone it from the row above
# (row 9) so it matches the rest of the table (border + base style).
$ws.Range("B9:D9").Copy()
$ws.Range("B10:D10").PasteSpecial(-4122)  # xlPasteFormats

# Populate the new "Step 2" row with the new model entry.
$ws.Range("B10").Value2 = 2
$ws.Range("C10").Value2 = "Model: Decision Tree"

# Renumber the Step column (B) for the rows that shifted down, 3..12.
$ws.Range("B11").Value2 = 3
$ws.Range("B12").Value2 = 4
$ws.Range("B13").Value2 = 5
$ws.Range("B14").Value2 = 6
$ws.Range("B15").Value2 = 7
$ws.Range("B16").Value2 = 8
$ws.Range("B17").Value2 = 9
$ws.Range("B18").Value2 = 10
$ws.Range("B19").Value2 = 11
$ws.Range("B20").Value2 = 12

# Restore the selection to C11 (matches the authored workbook state).
$ws.Range("C11").Select()
